$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.781.68'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.626.92'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.19'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5116'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.61%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2568'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06333'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.44'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07786'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.245'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.627.55'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.849.09'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5529'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.58'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₅7489'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.795.33'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.426'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '194.44'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.777'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.014'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.50%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.867'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.41'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1246'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.56'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.711'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.241'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04867'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.250'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.172'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.541'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.362'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8954'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5520'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.540'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.114.49'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01548'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.92%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.532'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7965'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.34'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.773.41'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -7.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4423'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.22%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.62'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.46%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.545'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.90%  '
